$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 265, shifting the existing rows 265:287 down to 266:288.
$ws.Rows.Item(265).Insert()

# Populate the newly inserted row 265 with the new weekly record.
$ws.Cells.Item(265, 1).Value = 4
$ws.Cells.Item(265, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(265, 3).Value = "Los Lagos"
$ws.Cells.Item(265, 4).Value = 44769
$ws.Cells.Item(265, 5).Value = 10
$ws.Cells.Item(265, 6).Value = 100112017
$ws.Cells.Item(265, 7).Value = "Apio"
$ws.Cells.Item(265, 8).Value = "Americana (o)"
$ws.Cells.Item(265, 9).Value = "Primera"
$ws.Cells.Item(265, 10).Value = 15
$ws.Cells.Item(265, 11).Value = 13000
$ws.Cells.Item(265, 12).Value = 13000
$ws.Cells.Item(265, 13).Value = 13000
$ws.Cells.Item(265, 14).Value = "`$/docena de matas"
$ws.Cells.Item(265, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(265, 16).Value = 2167
$ws.Cells.Item(265, 17).Value = 6
$ws.Cells.Item(265, 18).Value = "Hortaliza"
